$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 641-642 (shifts old rows 641.. down to 643..)
$ws.Range("A641:A642").EntireRow.Insert()

# Fill in the new row 641 (Larga vida / Primera, Arica y Parinacota)
$ws.Cells.Item(641,1).Value = 7
$ws.Cells.Item(641,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(641,3).Value = "Ñuble"
$ws.Cells.Item(641,4).Value = 45077
$ws.Cells.Item(641,5).Value = 16
$ws.Cells.Item(641,6).Value = 100112020
$ws.Cells.Item(641,7).Value = "Tomate"
$ws.Cells.Item(641,8).Value = "Larga vida"
$ws.Cells.Item(641,9).Value = "Primera"
$ws.Cells.Item(641,10).Value = 500
$ws.Cells.Item(641,11).Value = 15000
$ws.Cells.Item(641,12).Value = 16000
$ws.Cells.Item(641,13).Value = 15500
$ws.Cells.Item(641,14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(641,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(641,16).Value = 861
$ws.Cells.Item(641,17).Value = 18
$ws.Cells.Item(641,18).Value = "Hortaliza"

# Fill in the new row 642 (Larga vida / Segunda, Arica y Parinacota)
$ws.Cells.Item(642,1).Value = 7
$ws.Cells.Item(642,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(642,3).Value = "Ñuble"
$ws.Cells.Item(642,4).Value = 45077
$ws.Cells.Item(642,5).Value = 16
$ws.Cells.Item(642,6).Value = 100112020
$ws.Cells.Item(642,7).Value = "Tomate"
$ws.Cells.Item(642,8).Value = "Larga vida"
$ws.Cells.Item(642,9).Value = "Segunda"
$ws.Cells.Item(642,10).Value = 600
$ws.Cells.Item(642,11).Value = 12000
$ws.Cells.Item(642,12).Value = 13000
$ws.Cells.Item(642,13).Value = 12500
$ws.Cells.Item(642,14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(642,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(642,16).Value = 694
$ws.Cells.Item(642,17).Value = 18
$ws.Cells.Item(642,18).Value = "Hortaliza"
